$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume/Change% (E) columns per latest scrape.
# Force text format on Price cells so values like "1.00" / "0.0415" are
# preserved exactly as literal text rather than being coerced to numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.350.67"
$ws.Range("E2").Value = "  +0.41%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.416.91"
$ws.Range("E3").Value = "  +1.22%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.76"
$ws.Range("E5").Value = "  -0.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "176.17"
$ws.Range("E6").Value = "  -2.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.410.39"
$ws.Range("E8").Value = "  +1.09%  "
$ws.Range("E9").Value = "  -0.64%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.197"
$ws.Range("E10").Value = "  +0.60%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.580"
$ws.Range("E11").Value = "  -1.28%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "48.56"
$ws.Range("E12").Value = "  -0.09%  "
$ws.Range("E13").Value = "  -2.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "689.11"
$ws.Range("E14").Value = "  +0.49%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.967.12"
$ws.Range("E15").Value = "  +0.99%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.59"
$ws.Range("E16").Value = "  -0.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.437.28"
$ws.Range("E17").Value = "  +0.38%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.417.94"
$ws.Range("E18").Value = "  +0.76%  "
$ws.Range("E19").Value = "  +0.82%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.58"
$ws.Range("E20").Value = "  -0.59%  "
$ws.Range("E21").Value = "  -0.27%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.891"
$ws.Range("E22").Value = "  -0.86%  "
$ws.Range("E23").Value = "  +0.20%  "
$ws.Range("E24").Value = "  -1.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "101.06"
$ws.Range("E25").Value = "  -3.08%  "
$ws.Range("E26").Value = "  -0.99%  "
$ws.Range("E27").Value = "  -2.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.51"
$ws.Range("E28").Value = "  -0.74%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.31"
$ws.Range("E29").Value = "  -2.95%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.70"
$ws.Range("E30").Value = "  +0.37%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.01"
$ws.Range("E31").Value = "  +0.67%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "576.21"
$ws.Range("E32").Value = "  +3.76%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.67"
$ws.Range("E33").Value = "  +0.54%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "10.96"
$ws.Range("E34").Value = "  -1.95%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "58.21"
$ws.Range("E35").Value = "  +0.50%  "
$ws.Range("E36").Value = "  -3.30%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.559.98"
$ws.Range("E38").Value = "  -3.75%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.138"
$ws.Range("E39").Value = "  -1.56%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "34.72"
$ws.Range("E40").Value = "  -0.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0₃0724"
$ws.Range("E41").Value = "  +3.18%  "
$ws.Range("E42").Value = "  +0.53%  "
$ws.Range("E43").Value = "  -0.91%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.330"
$ws.Range("E44").Value = "  -2.26%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0415"
$ws.Range("E46").Value = "  +4.32%  "
$ws.Range("E47").Value = "  -0.04%  "
$ws.Range("E48").Value = "  -1.00%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.999"
$ws.Range("E49").Value = "  -0.31%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "132.44"
$ws.Range("E50").Value = "  -0.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.63"
$ws.Range("E51").Value = "  +2.51%  "
